# Apply updated crypto price/volume figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.743.16"
$ws.Range("E2").Value = "  -0.52%  "

$ws.Range("D3").Value = "1.596.14"
$ws.Range("E3").Value = "  -1.70%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "208.57"
$ws.Range("E5").Value = "  -1.35%  "

$ws.Range("D6").Value = "0.503"
$ws.Range("E6").Value = "  -2.36%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "22.39"
$ws.Range("E8").Value = "  -2.51%  "

$ws.Range("E9").Value = "  -1.55%  "

$ws.Range("E10").Value = "  -1.58%  "

$ws.Range("E11").Value = "  -1.40%  "

$ws.Range("D12").Value = "1.821.92"
$ws.Range("E12").Value = "  -1.86%  "

$ws.Range("D13").Value = "1.622.68"
$ws.Range("E13").Value = "  -0.39%  "

$ws.Range("E14").Value = "  -3.17%  "

$ws.Range("D15").Value = "0.535"
$ws.Range("E15").Value = "  -3.15%  "

$ws.Range("D16").Value = "27.728.06"
$ws.Range("E16").Value = "  -0.63%  "

$ws.Range("D17").Value = "63.57"
$ws.Range("E17").Value = "  -1.46%  "

$ws.Range("E18").Value = "  -3.16%  "

$ws.Range("D19").Value = "0.0₃0698"
$ws.Range("E19").Value = "  -2.27%  "

$ws.Range("D20").Value = "7.39"
$ws.Range("E20").Value = "  -2.72%  "

$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("E22").Value = "  -3.75%  "

$ws.Range("D23").Value = "9.75"
$ws.Range("E23").Value = "  -1.84%  "

$ws.Range("D24").Value = "1.99"
$ws.Range("E24").Value = "  -3.56%  "

$ws.Range("D25").Value = "154.03"
$ws.Range("E25").Value = "  -0.18%  "

$ws.Range("E26").Value = "  -1.61%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("E28").Value = "  -1.12%  "

$ws.Range("E29").Value = "  -4.55%  "

$ws.Range("D30").Value = "1.16"
$ws.Range("E30").Value = "  -1.25%  "

$ws.Range("E31").Value = "  -1.52%  "

$ws.Range("E32").Value = "  -4.50%  "

$ws.Range("D33").Value = "1.378.18"
$ws.Range("E33").Value = "  -2.40%  "

$ws.Range("E34").Value = "  -3.71%  "

$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  -3.41%  "

$ws.Range("D36").Value = "0.973"
$ws.Range("E36").Value = "  -2.38%  "

$ws.Range("E37").Value = "  -0.30%  "

$ws.Range("E38").Value = "  -0.72%  "

$ws.Range("E39").Value = "  -2.59%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.830"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.01%  "

$ws.Range("E41").Value = "  +0.07%  "

$ws.Range("D42").Value = "0.974"
$ws.Range("E42").Value = "  -2.85%  "

$ws.Range("D43").Value = "64.69"
$ws.Range("E43").Value = "  -0.84%  "

$ws.Range("E44").Value = "  +2.70%  "

$ws.Range("D45").Value = "5.22"
$ws.Range("E45").Value = "  -3.05%  "

$ws.Range("E46").Value = "  -4.04%  "

$ws.Range("D47").Value = "1.732.29"
$ws.Range("E47").Value = "  -2.01%  "

$ws.Range("D48").Value = "86.96"
$ws.Range("E48").Value = "  -2.06%  "

$ws.Range("D49").Value = "0.0₆0100"
$ws.Range("E49").Value = "  -1.49%  "

$ws.Range("D50").Value = "0.0967"
$ws.Range("E50").Value = "  -3.63%  "

$ws.Range("E51").Value = "  -1.19%  "
